$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 46.85851166666667
$ws.Range("H2").Value = 140.575535
$ws.Range("I2").Value = 0.1419057303676978
$ws.Range("J2").Value = 0.1419057303676978
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7328106666666666
$ws.Range("N2").Value = 2.198432
$ws.Range("O2").Value = 0.4072614640191846
$ws.Range("P2").Value = 0.4072614640191846
$ws.Range("Q2").Value = 34.33841717345777
$ws.Range("R2").Value = 309.04575456112
$ws.Range("S2").Value = 0.05779273550226026
$ws.Range("T2").Value = 0.05779273550226027

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 46.85851166666667
$ws.Range("H3").Value = 140.575535
$ws.Range("I3").Value = 0.1419057303676978
$ws.Range("J3").Value = 0.1419057303676978
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6260680000000001
$ws.Range("N3").Value = 1.878204
$ws.Range("O3").Value = 0.3479389450147599
$ws.Range("P3").Value = 0.3479389450147599
$ws.Range("Q3").Value = 29.33661468212667
$ws.Range("R3").Value = 264.02953213914
$ws.Range("S3").Value = 0.04937453011568575
$ws.Range("T3").Value = 0.04937453011568576

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 46.85851166666667
$ws.Range("H4").Value = 140.575535
$ws.Range("I4").Value = 0.1419057303676978
$ws.Range("J4").Value = 0.1419057303676978
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.440483
$ws.Range("N4").Value = 1.321449
$ws.Range("O4").Value = 0.2447995909660556
$ws.Range("P4").Value = 0.2447995909660555
$ws.Range("Q4").Value = 20.64037779446833
$ws.Range("R4").Value = 185.763400150215
$ws.Range("S4").Value = 0.03473846474975179
$ws.Range("T4").Value = 0.03473846474975179

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 283.1772156666667
$ws.Range("H5").Value = 849.531647
$ws.Range("I5").Value = 0.8575703363889615
$ws.Range("J5").Value = 0.8575703363889616
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7328106666666666
$ws.Range("N5").Value = 2.198432
$ws.Range("O5").Value = 0.4072614640191846
$ws.Range("P5").Value = 0.4072614640191846
$ws.Range("Q5").Value = 207.5152841975004
$ws.Range("R5").Value = 1867.637557777504
$ws.Range("S5").Value = 0.349255350697193
$ws.Range("T5").Value = 0.3492553506971931

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 283.1772156666667
$ws.Range("H6").Value = 849.531647
$ws.Range("I6").Value = 0.8575703363889615
$ws.Range("J6").Value = 0.8575703363889616
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6260680000000001
$ws.Range("N6").Value = 1.878204
$ws.Range("O6").Value = 0.3479389450147599
$ws.Range("P6").Value = 0.3479389450147599
$ws.Range("Q6").Value = 177.2881930579987
$ws.Range("R6").Value = 1595.593737521988
$ws.Range("S6").Value = 0.2983821181191281
$ws.Range("T6").Value = 0.2983821181191281

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 283.1772156666667
$ws.Range("H7").Value = 849.531647
$ws.Range("I7").Value = 0.8575703363889615
$ws.Range("J7").Value = 0.8575703363889616
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.440483
$ws.Range("N7").Value = 1.321449
$ws.Range("O7").Value = 0.2447995909660556
$ws.Range("P7").Value = 0.2447995909660555
$ws.Range("Q7").Value = 124.7347494885003
$ws.Range("R7").Value = 1122.612745396503
$ws.Range("S7").Value = 0.2099328675726405
$ws.Range("T7").Value = 0.2099328675726405

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.1730073333333333
$ws.Range("H8").Value = 0.519022
$ws.Range("I8").Value = 0.000523933243340694
$ws.Range("J8").Value = 0.000523933243340694
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7328106666666666
$ws.Range("N8").Value = 2.198432
$ws.Range("O8").Value = 0.4072614640191846
$ws.Range("P8").Value = 0.4072614640191846
$ws.Range("Q8").Value = 0.1267816192782222
$ws.Range("R8").Value = 1.141034573504
$ws.Range("S8").Value = 0.0002133778197312507
$ws.Range("T8").Value = 0.0002133778197312507

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.1730073333333333
$ws.Range("H9").Value = 0.519022
$ws.Range("I9").Value = 0.000523933243340694
$ws.Range("J9").Value = 0.000523933243340694
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6260680000000001
$ws.Range("N9").Value = 1.878204
$ws.Range("O9").Value = 0.3479389450147599
$ws.Range("P9").Value = 0.3479389450147599
$ws.Range("Q9").Value = 0.1083143551653333
$ws.Range("R9").Value = 0.974829196488
$ws.Range("S9").Value = 0.0001822967799461226
$ws.Range("T9").Value = 0.0001822967799461226

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1730073333333333
$ws.Range("H10").Value = 0.519022
$ws.Range("I10").Value = 0.000523933243340694
$ws.Range("J10").Value = 0.000523933243340694
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.440483
$ws.Range("N10").Value = 1.321449
$ws.Range("O10").Value = 0.2447995909660556
$ws.Range("P10").Value = 0.2447995909660555
$ws.Range("Q10").Value = 0.07620678920866665
$ws.Range("R10").Value = 0.685861102878
$ws.Range("S10").Value = 0.0001282586436633207
$ws.Range("T10").Value = 0.0001282586436633207
